$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 content ---
$ws.Range("A2").Value = 595
$ws.Range("B2").Value = "Big Countries"
$ws.Range("C2").Value = "SELECT"
$ws.Range("F2").Value = "Easy"
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = "✅"

# --- Row 3 content ---
$ws.Range("H3").Value = "❌"

# I2 set after H3 so new shared-string indices line up in the expected order
$ws.Range("I2").Value = "Given 2 sol and didn’t see solutions"

# --- Row height ---
$ws.Rows.Item(2).RowHeight = 28.5

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(9).ColumnWidth = 22.333333333333332

# --- Formatting ---
# H2 / H3: red font color for the check/cross marks
$ws.Range("H2").Font.Color = 255
$ws.Range("H3").Font.Color = 255

# I2 and the whole column I: wrap text
$ws.Range("I2").WrapText = $true
$ws.Columns.Item(9).WrapText = $true
# remove the incidental blank formatted cell the column-wide wrap created
$ws.Range("I3").Clear()

# A2: green fill
$ws.Range("A2").Interior.Color = 5296274

# --- Selection / view ---
$ws.Range("I3").Select()

# --- Page setup (adds <pageSetup orientation="portrait".../>) ---
$ws.PageSetup.Orientation = 1
